$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$conv = $wb.Worksheets.Item("CONVERTION")

$xlPasteFormats = -4122
$xlPasteValues  = -4163

$earnedFormula = '=IF(ISBLANK(Table1[[#This Row],[EARNED]]),"",Table1[[#This Row],[EARNED]])'

function Shift-Row($r) {
    $src = $ws.Range("A" + $r + ":K" + $r)
    $dst = $ws.Range("A" + ($r + 1) + ":K" + ($r + 1))
    $dst.ClearContents()
    $src.Copy()
    $dst.PasteSpecial($xlPasteFormats)
    $src.Copy()
    $dst.PasteSpecial($xlPasteValues)
}

# 1) Push rows 160..228 down to 161..229 (bottom-up so sources aren't clobbered
#    before they're read), opening up a blank row at 160 for the new entry.
for ($r = 228; $r -ge 160; $r--) {
    Shift-Row $r
}

# Re-establish the per-row "EARNED " calculated column formula (PasteValues
# flattens formulas to their cached results, so every shifted row needs it
# restored) for every data row from the new blank row through the new last row.
for ($r = 160; $r -le 229; $r++) {
    $ws.Cells.Item($r, 7).Formula = $earnedFormula
}

# 2) Grow the Table1 ListObject / sheet dimension to match the extra row.
$tbl = $ws.ListObjects.Item(1)
$tbl.Resize($ws.Range("A8:K229"))

# 3) Populate the newly opened row 160 ("UT(0-0-3)" leave usage entry).
$ws.Range("A160").ClearContents()
$ws.Range("B160").Value = "UT(0-0-3)"
$ws.Range("C160").ClearContents()
$ws.Range("D160").Value = 0.006
$ws.Range("E160").ClearContents()
$ws.Range("F160").ClearContents()
$ws.Range("H160").ClearContents()
$ws.Range("K160").ClearContents()

# 4) Record the new "UT(1-0-5)" leave usage on row 158.
$ws.Range("B158").Value = "UT(1-0-5)"
$ws.Range("D158").Value = 1.01

# 5) Update the CONVERTION sheet late-calculator inputs (1 day, 0 hours, 5
#    minutes) that back the 1.01 figure entered above.
$conv.Range("D3").Value = 1
$conv.Range("F3").Value = 5

# 6) Cosmetic: leave the selection where the user's edit session ended up.
$ws.Activate()
$ws.Range("F167").Select()
